$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-Text "D2" '43.361.84'
Set-Text "E2" '  +2.27%  '
Set-Text "D3" '2.392.39'
Set-Text "E3" '  +7.78%  '
Set-Text "E4" '  -0.25%  '
Set-Text "D5" '323.33'
Set-Text "E5" '  +11.45%  '
Set-Text "D6" '105.72'
Set-Text "E6" '  -4.90%  '
Set-Text "D7" '0.655'
Set-Text "E7" '  +4.67%  '
Set-Text "E8" '  -0.13%  '
Set-Text "D9" '0.653'
Set-Text "E9" '  +9.33%  '
Set-Text "D10" '41.85'
Set-Text "E10" '  -3.98%  '
Set-Text "D11" '0.0945'
Set-Text "E11" '  +3.83%  '
Set-Text "D12" '8.61'
Set-Text "E12" '  +0.06%  '
Set-Text "D13" '17.45'
Set-Text "E13" '  +17.41%  '
Set-Text "E14" '  +0.72%  '
Set-Text "E15" '  +2.70%  '
Set-Text "D16" '2.752.52'
Set-Text "E16" '  +7.71%  '
Set-Text "D17" '2.387.83'
Set-Text "E17" '  +7.10%  '
Set-Text "D18" '43.350.77'
Set-Text "E18" '  +2.34%  '
Set-Text "E19" '  +4.31%  '
Set-Text "D20" '7.41'
Set-Text "E20" '  +4.34%  '
Set-Text "D21" '76.24'
Set-Text "E21" '  +3.90%  '
Set-Text "B22" 'BitcoinCash'
Set-Text "C22" 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-Text "D22" '269.82'
Set-Text "E22" '  +14.85%  '
Set-Text "B23" 'PancakeSwap'
Set-Text "C23" 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-Text "D23" '3.44'
Set-Text "E23" '  +3.73%  '
Set-Text "E24" '  +2.11%  '
Set-Text "D25" '9.91'
Set-Text "E25" '  +11.88%  '
Set-Text "D26" '11.85'
Set-Text "E26" '  +4.31%  '
Set-Text "D27" '0.999'
Set-Text "E27" '  -0.05%  '
Set-Text "E28" '  +7.78%  '
Set-Text "D29" '177.18'
Set-Text "E29" '  +2.24%  '
Set-Text "B30" 'Toncoin'
Set-Text "C30" 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-Text "D30" '2.20'
Set-Text "E30" '  -0.10%  '
Set-Text "D31" '37.73'
Set-Text "E31" '  +0.83%  '
Set-Text "B32" 'WEMIXToken'
Set-Text "C32" 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-Text "D32" '3.21'
Set-Text "E32" '  +2.86%  '
Set-Text "D33" '0.0930'
Set-Text "E33" '  +6.42%  '
Set-Text "D34" '5.91'
Set-Text "E34" '  +5.77%  '
Set-Text "E35" '  +6.45%  '
Set-Text "D36" '4.88'
Set-Text "E36" '  -1.90%  '
Set-Text "D37" '4.10'
Set-Text "E37" '  -1.51%  '
Set-Text "D38" '0.0369'
Set-Text "E38" '  -2.19%  '
Set-Text "E39" '  +5.13%  '
Set-Text "E40" '  +18.81%  '
Set-Text "D41" '1.60'
Set-Text "E41" '  +22.47%  '
Set-Text "D42" '126.80'
Set-Text "E42" '  +25.54%  '
Set-Text "E43" '  +2.02%  '
Set-Text "D44" '69.52'
Set-Text "E44" '  -2.58%  '
Set-Text "E45" '  +0.09%  '
Set-Text "D46" '12.58'
Set-Text "E46" '  +2.37%  '
Set-Text "D47" '9.61'
Set-Text "E47" '  +14.57%  '
Set-Text "E48" '  +6.22%  '
Set-Text "D49" '88.22'
Set-Text "E49" '  +60.63%  '
Set-Text "E50" '  +3.00%  '
Set-Text "D51" '1.601.22'
Set-Text "E51" '  +12.60%  '
